$wb = $excel.ActiveWorkbook

# --- Add the new worksheet "Лист2" positioned right after "Кровля" ---
$ws2 = $wb.Worksheets.Add()
$ws2.Name = "Лист2"
$ws2.Move($null, $wb.Worksheets.Item("Кровля")) | Out-Null

# Re-fetch the sheet reference post-move (resolving fresh rather than
# reusing a pre-move object) and populate its two cells.
$ws2 = $wb.Worksheets.Item("Лист2")
$ws2.Range("A1").Value = "asaddd"
$ws2.Range("A2").Value = "ad"
$ws2.Range("C6").Select() | Out-Null

# --- Tidy up a stale "apply fill" flag on the "ФОТО" header cell (F1) ---
# Nudging the horizontal alignment away and back forces the style engine to
# rebuild the cell's xf record, which drops the leftover applyFill="1"
# attribute that no longer reflects any actual fill (fillId stays 0 / "no
# fill") while leaving the left alignment, font, border, etc. untouched.
$ws1 = $wb.Worksheets.Item("Кровля")
$f1 = $ws1.Range("F1")
$f1.HorizontalAlignment = -4152
$f1.HorizontalAlignment = -4131

# --- Update the remembered selection on "Кровля" ---
$ws1.Range("J29").Select() | Out-Null
